$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Sheet "汽車" (car) previously had its header row (row 1) wrongly
# duplicating the data row. Replace it with the proper column labels
# used by every other sheet in this workbook, and extend the table with
# the common metadata columns (property_category .. index).
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Fill in the new metadata columns (H:N) on the data row (row 2) to
# match the values already used by the other property sheets.
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-22"
$ws.Range("K2").Value = "馬文君"
$ws.Range("L2").Value = 1724
$ws.Range("M2").Value = "tmp99351"
$ws.Range("N2").Value = 41

# Match the formatting of the new cells to their neighbouring columns
# (bold/bordered header style for row 1, plain style for row 2).
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
